# "Generate Report for Archive"
# - Flip the handoff/translation status shown on the Overview sheet (both
#   language columns) and on each per-language detail sheet from
#   "Ready for handoff" to "In Translation".
# - Narrow the per-language status columns (Overview!E:F and the "Status"
#   column on each language sheet) to their tighter archived-report width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
